$d = $word.ActiveDocument

$d.Content.Find.Execute("Know how to code in Java,", $true, $false, $false, $false, $false, $true, 1, $false, "Know how to code in Javascript,", 2)
